$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B15: ripple frequency factor now pulled from the simulation table (F17)
# instead of the hard-coded 100kHz duty-cycle figure.
$ws.Range("B15").Formula = "=F17"

# B17: H1 inductance now pulled from the simulation table (F18) instead of
# the previously hard-coded 2600 value.
$ws.Range("B17").Formula = "=F18"

# Update the current selection / view to reflect where the user ended up
# after testing the simulation with the new input voltages.
$ws.Range("H10").Select()
